$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to Text format before writing so that numeric-looking
# strings (e.g. "3.74", "206.34") are preserved verbatim as text instead of
# being coerced into floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @(
    @("D2", "26.919.28"),
    @("D3", "1.550.52"),
    @("E3", "  -0.11%  "),
    @("E4", "  -0.37%  "),
    @("D5", "206.34"),
    @("E5", "  +0.08%  "),
    @("E6", "  +0.93%  "),
    @("E7", "  -0.36%  "),
    @("D8", "22.08"),
    @("E8", "  +2.64%  "),
    @("E10", "  +0.87%  "),
    @("E11", "  -0.27%  "),
    @("D12", "1.771.27"),
    @("E12", "  -0.08%  "),
    @("D13", "1.549.33"),
    @("E13", "  +0.34%  "),
    @("D14", "3.74"),
    @("E14", "  +0.91%  "),
    @("D15", "0.519"),
    @("D16", "26.898.45"),
    @("E16", "  -0.01%  "),
    @("D17", "61.66"),
    @("E17", "  +0.09%  "),
    @("D18", "217.16"),
    @("E18", "  +1.53%  "),
    @("D19", "0.0₃0699"),
    @("E19", "  +1.98%  "),
    @("D20", "7.27"),
    @("E20", "  +0.50%  "),
    @("E21", "  -0.42%  "),
    @("D22", "4.06"),
    @("E22", "  +0.40%  "),
    @("E24", "  -1.17%  "),
    @("D25", "153.90"),
    @("E25", "  +0.35%  "),
    @("E26", "  -0.39%  "),
    @("D27", "14.97"),
    @("E27", "  +0.68%  "),
    @("E28", "  +0.80%  "),
    @("E29", "  -0.37%  "),
    @("E30", "  +2.01%  "),
    @("E31", "  -0.21%  "),
    @("E32", "  -0.20%  "),
    @("B33", "Maker"),
    @("C33", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"),
    @("D33", "1.414.25"),
    @("E33", "  +3.27%  "),
    @("B34", "InternetComputer(DFINITY)"),
    @("C34", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"),
    @("D34", "3.10"),
    @("E34", "  +4.86%  "),
    @("E35", "  +2.93%  "),
    @("E36", "  -0.65%  "),
    @("E37", "  +0.06%  "),
    @("D39", "0.527"),
    @("E39", "  +1.16%  "),
    @("E40", "  +0.06%  "),
    @("E41", "  -0.36%  "),
    @("D42", "5.70"),
    @("E42", "  +3.38%  "),
    @("E43", "  +3.01%  "),
    @("D44", "0.999"),
    @("E44", "  +1.41%  "),
    @("D45", "64.68"),
    @("E45", "  +1.76%  "),
    @("E46", "  +0.91%  "),
    @("D47", "1.685.38"),
    @("E47", "  -0.05%  "),
    @("D48", "87.51"),
    @("E48", "  +1.64%  "),
    @("D49", "0.0515"),
    @("E49", "  +1.72%  "),
    @("E50", "  +3.54%  "),
    @("D51", "0.0961"),
    @("E51", "  +0.91%  ")
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# Restore default (Normal) style on the range so no stray number-format
# style index is left attached to the cells.
$ws.Range("D2:E51").Style = "Normal"

Write-Host "Applied $($updates.Count) cell updates"
